# Overwritten functionality of export of Purchase Order.
# Now 2 different dataProviders are responsible for viewing the grid on
# HTML and creating the Excel file -- update the exported values so the
# sheet reflects the new (Excel specific) data provider's output.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# SKU values produced by the new Excel data provider
$ws.Range("A3").Value = "1007e2"
$ws.Range("A4").Value = "31903"

# Supplier field that used to be left blank now has a value
$ws.Range("H4").Value = "Washka Pashka"

# Page setup for printing the generated workbook
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Cursor/selection ends up on P13 after the export routine runs
$ws.Range("P13").Select() | Out-Null
